$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day's price row is published on top of the table (latest date first),
# so insert a fresh row 2 and let everything else cascade down by one.
$ws.Rows("2:2").Insert()

# Force the date-like text columns (A: Date, E: Circular Date) to be entered
# as literal text rather than being auto-parsed into real date serials, so
# the stored cell stays a plain string like the rest of the sheet.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

# New top row: 11-11-2025, same circular/price info as the (now shifted down)
# row that used to be on top (row 3 after the insert == old row 2).
$ws.Range("A2").Value = "11-11-2025"
$ws.Range("B2").Value = $ws.Range("B3").Value()
$ws.Range("C2").Value = $ws.Range("C3").Value()
$ws.Range("D2").Value = $ws.Range("D3").Value()
$ws.Range("E2").Value = $ws.Range("E3").Value()
$ws.Range("F2").Value = $ws.Range("F3").Value()

# Re-apply the original cell formatting (style indexes) from row 3 onto the
# freshly inserted row 2, undoing the NumberFormat tweak above and the
# bold/header-ish formatting that Insert() pulled down from row 1.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The row-insert already shifted/duplicated the F2:F97 hyperlinks correctly
# (F2 now reuses the same hyperlink as F3). The very last row (old F97,
# now F98) has no hyperlink yet because nothing existed below it to copy
# from, so add it explicitly, pointing at the same circular PDF link shown
# in F98's text.
$ws.Hyperlinks.Add($ws.Range("F98"), $ws.Range("F98").Value())
